$wb = $excel.ActiveWorkbook

# Overview sheet: row 4 (e013cfbb...) "Latest HO Xliff Generate Date".
# This text is identical to (and shares the same backing shared-string as)
# de-de!H4 "Correspond Handoff Datetime" for the same file, so both cells
# must be updated together to keep sharing that string.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-09-04 16:50:05"

# zh-cn sheet: row 4 (e013cfbb...) Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-09-04 16:49:57"
$wsZhCn.Range("K4").Value = "2016-09-04 16:50:34"

# de-de sheet: row 4 (e013cfbb...) Correspond Handoff Datetime / Correspond Handback DateTime
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H4").Value = "2016-09-04 16:50:05"
$wsDeDe.Range("K4").Value = "2016-09-04 16:50:42"
